$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.392.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +6.72%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.816.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +6.64%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "343.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3857"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.80%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "50.50"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3541"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.244"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07815"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.75"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +14.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.002"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.683"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.272"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.814.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.95%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001137"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06776"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "86.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.69%  "
$ws.Range("E20").Value = "  +0.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +11.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.592"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +9.46%  "
$ws.Range("E23").Value = "  +2.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "27.383.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.468"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.763"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.95%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +15.35%  "
$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.520"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +20.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "154.10"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.017.86"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "137.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.465"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +9.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.125"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "13.98"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08846"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.727"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.55%  "
$ws.Range("B37").Value = "TheSandbox"
$ws.Range("C37").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.7226"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +18.72%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.683"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06596"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02436"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2281"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.098"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.270"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6760"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +15.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.974"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.207"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "134.13"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07351"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.40"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.70%  "
